# Add "Other1" (Y) and "Total1" (Z) computed columns to the KO comparison sheet.
# This inserts two new columns before the existing "Other" column, pushing the
# previous Other/Total/Assembler/Season columns two slots to the right
# (Y->AA, Z->AB, AA->AC, AB->AD), and fills the new columns with per-row
# subtotal formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at Y:Z - everything from the old Y onward shifts right.
$ws.Columns("Y:Z").Insert()

# New column headers.
$ws.Range("Y1").Value = "Other1"
$ws.Range("Z1").Value = "Total1"

# Per-row formulas for the two new columns, and refresh the (now shifted)
# grand-total column AB so it references the new AA ("Other") column instead
# of the old Y.
$lastRow = 19
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("Y$r").Formula = "=SUM(F$r+I$r+J$r+K$r+L$r+M$r+N$r+O$r+Q$r+R$r+S$r+T$r+V$r+W$r+X$r)"
    $ws.Range("Z$r").Formula = "=SUM(Y$r+D$r+E$r+G$r+H$r)"
    $ws.Range("AB$r").Formula = "=C$r+P$r+U$r+AA$r"
}

# Match the saved selection from the authored workbook.
$ws.Range("W28").Select()
